# "Generate Report for handoff" — the handoff transform failed for the
# e2e test file, so the localization-status report now points at a new
# transient filename and reflects the failure instead of a successful
# "Ready for handoff" / completed-handback state.

$wb = $excel.ActiveWorkbook

$oldName = "2c026530-0b0e-4853-a10c-c1820915f1f6.md"
$newName = "76f98c9c-c65b-44b5-9710-80dc295ad7b5.md"
$newUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/6189284bef85735aa3ff72c9f28a22fe9c1401ab/e2e/$newName"

# ---- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value2 = $newName
$ov.Range("B2").Value2 = "Handoff transform failed"
$ov.Range("C2").Value2 = "Handoff transform failed"

$ov.Cells.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), $newUrl, "", "", $newName) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6189284bef85735aa3ff72c9f28a22fe9c1401ab/.localization-config", "", "", ".localization-config") | Out-Null

# ---- per-language detail sheets -------------------------------------
$langSheets = @("zh-cn", "de-de")
foreach ($langName in $langSheets) {
    $ws = $wb.Worksheets.Item($langName)

    $ws.Range("A2").Value2 = $newName
    $ws.Range("B2").Value2 = "Handoff transform failed"
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value2 = "0001-01-01 00:00:00"
    $ws.Range("G2").Value2 = "0001-01-01 00:00:00"
    $ws.Range("H2").Value2 = "Ignored"

    $ws.Range("D3").Value2 = "0001-01-01 00:00:00"
    $ws.Range("G3").Value2 = "0001-01-01 00:00:00"
    $ws.Range("H3").Value2 = "Ignored"

    $ws.Cells.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $newUrl, "", "", $newName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6189284bef85735aa3ff72c9f28a22fe9c1401ab/.localization-config", "", "", ".localization-config") | Out-Null
}
